$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 543, shifting all existing rows (543-571) down by one (to 544-572)
$ws.Rows.Item(543).Insert()

# Fill the new row 543 with the latest weekly price entry
$ws.Range("A543").Value = 4
$ws.Range("B543").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C543").Value = "Los Lagos"
$ws.Range("D543").Value = 45267
$ws.Range("E543").Value = 10
$ws.Range("F543").Value = 100114014
$ws.Range("G543").Value = "Betarraga"
$ws.Range("H543").Value = "Sin especificar"
$ws.Range("I543").Value = "Primera"
$ws.Range("J543").Value = 500
$ws.Range("K543").Value = 1000
$ws.Range("L543").Value = 1100
$ws.Range("M543").Value = 1050
$ws.Range("N543").Value = "$/paquete 5 unidades"
$ws.Range("O543").Value = "Región Metropolitana"
$ws.Range("P543").Value = 210
$ws.Range("Q543").Value = 5
$ws.Range("R543").Value = "Hortaliza"
